$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Overview sheet: handback status text changed for both locale rows
# -----------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# Hyperlink font color used for the "Latest Target File" cells below.
# FF6495ED (ARGB) == BGR 0xED9564 == 15570276 decimal, matching the
# workbook's existing HyperLink cell style.
$hyperlinkColor = 15570276

# -----------------------------------------------------------------
# zh-cn sheet: fill in "Latest Target File" / "Latest Handback File"
# / "Latest Handback DateTime" now that the handback is complete
# -----------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7aa080c7600c6d98bd99ec17933ae8710caf0993/e2e/0d44ddd4-919f-4ae6-bb1d-8f7545878d1d.md",
    "",
    "",
    "0d44ddd4-919f-4ae6-bb1d-8f7545878d1d.md"
) | Out-Null
$wsZhCn.Range("J2").Value = "0d44ddd4-919f-4ae6-bb1d-8f7545878d1d.87b437dd0e36bf5f8c10db809d4443ed0e6f4e3e.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-15 18:26:54"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7aa080c7600c6d98bd99ec17933ae8710caf0993/e2e/42069519-92a2-4de6-8c14-aa7fad3efa75.md",
    "",
    "",
    "42069519-92a2-4de6-8c14-aa7fad3efa75.md"
) | Out-Null
$wsZhCn.Range("J3").Value = "42069519-92a2-4de6-8c14-aa7fad3efa75.6a6987cbe0afe10fe1723301d814d9e261250861.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-15 18:26:54"

$wsZhCn.Range("I2:I3").Font.Underline = $true
$wsZhCn.Range("I2:I3").Font.Color = $hyperlinkColor

# -----------------------------------------------------------------
# de-de sheet: same handback fields, completed a little later
# -----------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7aa080c7600c6d98bd99ec17933ae8710caf0993/e2e/0d44ddd4-919f-4ae6-bb1d-8f7545878d1d.md",
    "",
    "",
    "0d44ddd4-919f-4ae6-bb1d-8f7545878d1d.md"
) | Out-Null
$wsDeDe.Range("J2").Value = "0d44ddd4-919f-4ae6-bb1d-8f7545878d1d.87b437dd0e36bf5f8c10db809d4443ed0e6f4e3e.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-15 18:27:03"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7aa080c7600c6d98bd99ec17933ae8710caf0993/e2e/42069519-92a2-4de6-8c14-aa7fad3efa75.md",
    "",
    "",
    "42069519-92a2-4de6-8c14-aa7fad3efa75.md"
) | Out-Null
$wsDeDe.Range("J3").Value = "42069519-92a2-4de6-8c14-aa7fad3efa75.6a6987cbe0afe10fe1723301d814d9e261250861.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-15 18:27:03"

$wsDeDe.Range("I2:I3").Font.Underline = $true
$wsDeDe.Range("I2:I3").Font.Color = $hyperlinkColor
